$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-12 Tuesday", "2024-11-13 Wednesday"),
    @("439×5=2195", "611×8=4888"),
    @("672×8=5376", "208×8=1664"),
    @("997×8=7976", "105×5=525"),
    @("648×3=1944", "272×3=816"),
    @("320×7=2240", "557×7=3899"),
    @("165×5=825", "717×4=2868"),
    @("302×8=2416", "973×4=3892"),
    @("889×8=7112", "851×3=2553"),
    @("836×7=5852", "835×4=3340"),
    @("159×9=1431", "134×8=1072"),
    @("469×8=3752", "433×4=1732"),
    @("617×2=1234", "904×3=2712"),
    @("769×2=1538", "748×8=5984"),
    @("397×3=1191", "369×8=2952"),
    @("572×6=3432", "641×7=4487"),
    @("399×9=3591", "566×5=2830"),
    @("926×2=1852", "102×6=612"),
    @("503×5=2515", "178×2=356"),
    @("868×9=7812", "105×4=420"),
    @("857×8=6856", "949×2=1898"),
    @("766×6=4596", "909×8=7272"),
    @("262×3=786", "989×6=5934"),
    @("692×3=2076", "714×5=3570"),
    @("334×7=2338", "572×9=5148"),
    @("771×2=1542", "219×7=1533")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
